# This workbook's sheet is protected; unprotect it so the cells can be edited,
# then re-protect (without a password, since the original legacy password hash
# cannot be reproduced) to leave the sheet in a protected state afterwards.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$wasProtected = $ws.ProtectContents
$ws.Unprotect()

# Update the "as of" date in the confidential disclaimer text (row 13, col A)
$disclaimer = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-05 for illustrative purposes only and are subject to change."
$ws.Range("A13").Value = $disclaimer
# Re-fit the row height back to automatic, since Excel otherwise stamps a
# fixed custom row height onto row 13 after the multi-line text is rewritten.
$ws.Rows.Item(13).AutoFit()

# Update the Weight (D) and Percent Change (E) columns for rows 2-9, and
# the Percent Change for row 10 (Total row), with the newly refreshed figures.
$ws.Range("D2").Value = 0.09722849652088532
$ws.Range("E2").Value = -0.01594341884440176

$ws.Range("D3").Value = 0.1073585921371751
$ws.Range("E3").Value = -0.005426680183672206

$ws.Range("D4").Value = 0.1194523592399912
$ws.Range("E4").Value = -0.000626664577784819

$ws.Range("D5").Value = 0.1402161906766053
$ws.Range("E5").Value = 0.001034949446700129

$ws.Range("D6").Value = 0.1361125408003148
$ws.Range("E6").Value = -0.0008321775312065682

$ws.Range("D7").Value = 0.1447194038699169
$ws.Range("E7").Value = 0.007622677465459526

$ws.Range("D8").Value = 0.1263713381577403
$ws.Range("E8").Value = 0.005975500448162752

$ws.Range("D9").Value = 0.1285410785973712
$ws.Range("E9").Value = -0.008183239588198443

$ws.Range("E10").Value = -0.001369366196217858

if ($wasProtected) {
    $ws.Protect()
}

$wb.Save()
